$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row3
$ws.Range("B3").Value = -50.75184683631875
$ws.Range("C3").Value = -40.75184683631875
$ws.Range("D3").Value = -30.75184683631874
$ws.Range("E3").Value = -25.75184683631874
$ws.Range("F3").Value = -25.75184683631874
$ws.Range("G3").Value = -25.75184683631874
$ws.Range("H3").Value = -25.75184683631874
$ws.Range("I3").Value = -25.75184683631874
$ws.Range("J3").Value = -20.75184683631874
$ws.Range("K3").Value = -10.75184683631874
$ws.Range("L3").Value = -0.7518468363187445
$ws.Range("M3").Value = 19.24815316368126
$ws.Range("N3").Value = 39.24815316368125
$ws.Range("O3").Value = 59.24815316368125

# row4
$ws.Range("F4").Value = 2.730995758697492
$ws.Range("G4").Value = 1.722369395643906
$ws.Range("H4").Value = 2.878921494922184
$ws.Range("I4").Value = 1.859543089110336

# row5
$ws.Range("B5").Value = 77.11499999999999
$ws.Range("C5").Value = 70.98999999999999
$ws.Range("D5").Value = 65.663
$ws.Range("E5").Value = 63.272
$ws.Range("F5").Value = 123.12
$ws.Range("G5").Value = 90.684
$ws.Range("H5").Value = 127.508
$ws.Range("I5").Value = 95.41
$ws.Range("J5").Value = 61.048
$ws.Range("K5").Value = 57.053
$ws.Range("L5").Value = 53.586
$ws.Range("M5").Value = 47.922
$ws.Range("N5").Value = 43.528
$ws.Range("O5").Value = 40.038

# row6
$ws.Range("B6").Value = 37.833
$ws.Range("C6").Value = 34.828
$ws.Range("D6").Value = 32.214
$ws.Range("E6").Value = 31.041
$ws.Range("F6").Value = 60.403
$ws.Range("G6").Value = 44.49
$ws.Range("H6").Value = 62.555
$ws.Range("I6").Value = 46.808
$ws.Range("J6").Value = 29.95
$ws.Range("K6").Value = 27.99
$ws.Range("L6").Value = 26.289
$ws.Range("M6").Value = 23.511
$ws.Range("N6").Value = 21.355
$ws.Range("O6").Value = 19.643

# row7
$ws.Range("B7").Value = 2350.487
$ws.Range("C7").Value = 2163.795
$ws.Range("D7").Value = 2001.426
$ws.Range("E7").Value = 1928.548
$ws.Range("F7").Value = 1374.126
$ws.Range("G7").Value = 1604.809
$ws.Range("H7").Value = 1349.977
$ws.Range("I7").Value = 1563.891
$ws.Range("J7").Value = 1860.76
$ws.Range("K7").Value = 1738.991
$ws.Range("L7").Value = 1633.316
$ws.Range("M7").Value = 1460.676
$ws.Range("N7").Value = 1326.745
$ws.Range("O7").Value = 1220.369

# row8
$ws.Range("B8").Value = 31.467
$ws.Range("C8").Value = 28.968
$ws.Range("D8").Value = 26.794
$ws.Range("E8").Value = 25.818
$ws.Range("F8").Value = 50.24
$ws.Range("G8").Value = 37.004
$ws.Range("H8").Value = 52.029
$ws.Range("I8").Value = 38.932
$ws.Range("J8").Value = 24.911
$ws.Range("K8").Value = 23.28
$ws.Range("L8").Value = 21.866
$ws.Range("M8").Value = 19.555
$ws.Range("N8").Value = 17.762
$ws.Range("O8").Value = 16.338

# row9
$ws.Range("B9").Value = 4.077
$ws.Range("C9").Value = 4.429
$ws.Range("D9").Value = 4.789
$ws.Range("E9").Value = 4.97
$ws.Range("F9").Value = 6.978
$ws.Range("G9").Value = 5.973
$ws.Range("H9").Value = 7.103
$ws.Range("I9").Value = 6.13
$ws.Range("J9").Value = 5.151
$ws.Range("K9").Value = 5.512
$ws.Range("L9").Value = 5.869
$ws.Range("M9").Value = 6.564
$ws.Range("N9").Value = 7.227
$ws.Range("O9").Value = 7.858

# row10
$ws.Range("B10").Value = 5.629
$ws.Range("C10").Value = 6.116
$ws.Range("D10").Value = 6.612
$ws.Range("E10").Value = 6.862
$ws.Range("F10").Value = 9.637
$ws.Range("G10").Value = 8.249000000000001
$ws.Range("H10").Value = 9.81
$ws.Range("I10").Value = 8.465
$ws.Range("J10").Value = 7.113
$ws.Range("K10").Value = 7.612
$ws.Range("L10").Value = 8.105
$ws.Range("M10").Value = 9.065
$ws.Range("N10").Value = 9.981999999999999
$ws.Range("O10").Value = 10.854

# row11
$ws.Range("B11").Value = 3.605
$ws.Range("C11").Value = 3.916
$ws.Range("D11").Value = 4.234
$ws.Range("E11").Value = 4.394
$ws.Range("F11").Value = 6.169
$ws.Range("G11").Value = 5.281
$ws.Range("H11").Value = 6.28
$ws.Range("I11").Value = 5.42
$ws.Range("J11").Value = 4.554
$ws.Range("K11").Value = 4.873
$ws.Range("L11").Value = 5.189
$ws.Range("M11").Value = 5.803
$ws.Range("N11").Value = 6.39
$ws.Range("O11").Value = 6.948

# row12
$ws.Range("B12").Value = 4.788
$ws.Range("C12").Value = 5.201
$ws.Range("D12").Value = 5.624
$ws.Range("E12").Value = 5.836
$ws.Range("F12").Value = 8.195
$ws.Range("G12").Value = 7.015
$ws.Range("H12").Value = 8.342000000000001
$ws.Range("I12").Value = 7.199
$ws.Range("J12").Value = 6.049
$ws.Range("K12").Value = 6.473
$ws.Range("L12").Value = 6.893
$ws.Range("M12").Value = 7.709
$ws.Range("N12").Value = 8.489000000000001
$ws.Range("O12").Value = 9.23

# row13
$ws.Range("B13").Value = 1.049
$ws.Range("C13").Value = 1.139
$ws.Range("D13").Value = 1.231
$ws.Range("E13").Value = 1.278
$ws.Range("F13").Value = 1.794
$ws.Range("G13").Value = 1.536
$ws.Range("H13").Value = 1.826
$ws.Range("I13").Value = 1.576
$ws.Range("J13").Value = 1.325
$ws.Range("K13").Value = 1.417
$ws.Range("L13").Value = 1.509
$ws.Range("M13").Value = 1.687
$ws.Range("N13").Value = 1.858
$ws.Range("O13").Value = 2.02
